$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Shopping List: Week of Sunday, February 05, 2022" "Shopping List: Week of Sunday, February 20, 2022"
Replace-Text "Veal Scallopini, Bacon" "Steak"
Replace-Text "Eggs, Lipton Envelope" "Lipton Envelope, Eggs"
Replace-Text "Parmigiano, Pecorino, Milk" "Parmigiano, Butter, Ricotta, Milk"
Replace-Text "Zucchini, Romaine Lettuce, Tomato, Asparagus, Peas, Onion" "Tomatoes, Onion, Garlic, Basil, Romaine Lettuce, Vegetable Medley, Asparagus, Peas, Broccoli Florets"
Replace-Text "Flour, White Bread, Baguette" "White Bread, Baguette"
Replace-Text "Spaghetti, Ditali" "Penne, Ditali"
Replace-Text "Pepper, Salt, Bread Crumbs" "Salt, Pepper"
Replace-Text "Olive Oil, White Vinegar, White Balsamic Vinegar" "Olive Oil, White Balsamic Vinegar"
Replace-Text "Menu: Week of Saturday, February 05, 2022" "Menu: Week of Sunday, February 20, 2022"
Replace-Text "Cacio e Pepe + Boiled Zucchini" "Pomodori + Salad"
Replace-Text "Breaded Cutlets + Salad" "American Steak + Vegetable Medley"
Replace-Text "Carbonara + Salad" "Ricotta and Broccoli"
